# Presupuesto Conceptos para revision
#
# Updates the "Hacker Etico" contractor rate from 2181.82 to 1363.63/1363.64
# (and the one place it is hard-coded a second time on "Flujo de Caja"),
# bumps the "Valor Actual Neto" discount-rate input on B12, re-points the
# custom percentage style at a 3-decimal format, and leaves the window
# focused on the first sheet.

$wb = $excel.ActiveWorkbook

$wsConceptos = $wb.Worksheets.Item("Presupuesto por conceptos")
$wsTareas    = $wb.Worksheets.Item("Presupuesto por tareas")
$wsFlujo     = $wb.Worksheets.Item("Flujo de Caja")
$wsVAN       = $wb.Worksheets.Item("Valor Actual Neto")

# ---------------------------------------------------------------------------
# 1. "Presupuesto por conceptos" - Hacker Etico rate: 2181.82 -> 1363.63
#    Every other value on this sheet (B3, B4, B5, B6, B22, B24) is a formula
#    that depends on B23, so it recalculates automatically.
# ---------------------------------------------------------------------------
$wsConceptos.Range("B23").Formula = "=1363.63*3.5"

# ---------------------------------------------------------------------------
# 2. "Flujo de Caja" - same rate change, spread across the row-14 "Hacker
#    Etico" line (quarterly pattern in columns B/I/P/Y). The formulas that
#    used to recompute this each quarter are flattened to literals/updated
#    constants, matching how Excel leaves a shared-formula block after only
#    part of it is touched by hand.
# ---------------------------------------------------------------------------
$wsFlujo.Range("B14").Value = 0
$wsFlujo.Range("I14").Formula = "=-1363.64*3.5/3"
$wsFlujo.Range("P14").Value = -1590.9133333333336
$wsFlujo.Range("Y14").Value = -1590.9133333333336
$wsFlujo.Range("B15").Value = 0

# ---------------------------------------------------------------------------
# 3. "Valor Actual Neto" - discount rate input B12: 10.57% -> 12.255%
#    (B5 keeps its value but picks up the same refreshed 3-decimal display).
# ---------------------------------------------------------------------------
$wsVAN.Range("B12").Value = 0.12255000000000001
$wsVAN.Range("B5").NumberFormat = "0.000%"
$wsVAN.Range("B12").NumberFormat = "0.000%"

# ---------------------------------------------------------------------------
# 4. View / window state bookkeeping.
# ---------------------------------------------------------------------------
$wsFlujo.Activate()
$wsFlujo.Range("B14").Select()

$wsVAN.Activate()
$wsVAN.Range("C7").Select()

$wsConceptos.Activate()
$wsConceptos.Range("E20").Select()

$excel.ActiveWindow.WindowState = -4140
